# "correccion en df de autos y motos"
# IND_VIA_MOTOS sheet: a new "EL PEÑON - SABADO" / "EL PEÑON - DOMINGO" pair of
# blocks is inserted right after "EL PEÑON - TIPICO"; every block from
# "GRANADA - TIPICO" onward shifts down by one 7-row slot, and a brand new
# "SAN FERNANDO (PARQUE DEL PERRO) - SABADO" block is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IND_VIA_MOTOS")

# --- Grow the formatted block template down to rows 52-70 -------------------
# Rows 52:56 / 59:63 / 66:70 don't exist yet in the source sheet, so clone the
# row formatting (styles, borders, fills, widths) from existing 7-row blocks
# before overwriting their contents below.
$ws.Range("A38:B42").Copy($ws.Range("A52:B56"))
$ws.Range("A45:B49").Copy($ws.Range("A59:B63"))
$ws.Range("A24:B28").Copy($ws.Range("A66:B70"))

# The section-title row of each template block (e.g. A38) has no sibling B
# cell in the source data; the Copy above still stamps a blank, formatted B
# cell next to the new title rows. Drop those blanks so the header rows only
# carry the A cell, matching every other section-title row on this sheet.
$ws.Cells.Item(52, 2).ClearContents()
$ws.Cells.Item(59, 2).ClearContents()
$ws.Cells.Item(66, 2).ClearContents()

# --- EL PEÑON - TIPICO (row 10) : Demanda total changes ---------------------
$ws.Cells.Item(13, 2).Value = 816

# --- EL PEÑON - SABADO (new content at the old GRANADA - TIPICO slot) -------
$ws.Cells.Item(17, 1).Value = "EL PEÑON - SABADO"
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(20, 2).Value = 3
$ws.Cells.Item(21, 2).Value = "0:00"

# --- EL PEÑON - DOMINGO (new content at the old GRANADA - SABADO slot) -----
$ws.Cells.Item(24, 1).Value = "EL PEÑON - DOMINGO"
$ws.Cells.Item(26, 2).Value = 3
$ws.Cells.Item(27, 2).Value = 4
$ws.Cells.Item(28, 2).Value = "0:00"

# --- GRANADA - TIPICO (content shifted from old GRANADA - DOMINGO slot) ----
$ws.Cells.Item(31, 1).Value = "GRANADA - TIPICO"
$ws.Cells.Item(33, 2).Value = 176
$ws.Cells.Item(34, 2).Value = 766
$ws.Cells.Item(35, 2).Value = "3:01"

# --- GRANADA - SABADO (content shifted from old SAN ANTONIO - TIPICO slot) -
$ws.Cells.Item(38, 1).Value = "GRANADA - SABADO"
$ws.Cells.Item(40, 2).Value = 109
$ws.Cells.Item(41, 2).Value = 441
$ws.Cells.Item(42, 2).Value = "1:47"

# --- GRANADA - DOMINGO (content shifted from old SAN FERNANDO TIPICO slot) -
$ws.Cells.Item(45, 1).Value = "GRANADA - DOMINGO"
$ws.Cells.Item(47, 2).Value = 157
$ws.Cells.Item(48, 2).Value = 479
$ws.Cells.Item(49, 2).Value = "1:48"

# --- SAN ANTONIO - TIPICO (new block) ---------------------------------------
$ws.Cells.Item(52, 1).Value = "SAN ANTONIO - TIPICO"
$ws.Cells.Item(53, 1).Value = "INDICADOR"
$ws.Cells.Item(53, 2).Value = "VALOR"
$ws.Cells.Item(54, 1).Value = "Ocupación Máxima"
$ws.Cells.Item(54, 2).Value = 106
$ws.Cells.Item(55, 1).Value = "Demanda total"
$ws.Cells.Item(55, 2).Value = 469
$ws.Cells.Item(56, 1).Value = "Duración Media (Dm)"
$ws.Cells.Item(56, 2).Value = "3:07"

# --- SAN FERNANDO (PARQUE DEL PERRO) - TIPICO (new block) ------------------
$ws.Cells.Item(59, 1).Value = "SAN FERNANDO (PARQUE DEL PERRO) - TIPICO"
$ws.Cells.Item(60, 1).Value = "INDICADOR"
$ws.Cells.Item(60, 2).Value = "VALOR"
$ws.Cells.Item(61, 1).Value = "Ocupación Máxima"
$ws.Cells.Item(61, 2).Value = 184
$ws.Cells.Item(62, 1).Value = "Demanda total"
$ws.Cells.Item(62, 2).Value = 907
$ws.Cells.Item(63, 1).Value = "Duración Media (Dm)"
$ws.Cells.Item(63, 2).Value = "2:33"

# --- SAN FERNANDO (PARQUE DEL PERRO) - SABADO (new block) ------------------
$ws.Cells.Item(66, 1).Value = "SAN FERNANDO (PARQUE DEL PERRO) - SABADO"
$ws.Cells.Item(67, 1).Value = "INDICADOR"
$ws.Cells.Item(67, 2).Value = "VALOR"
$ws.Cells.Item(68, 1).Value = "Ocupación Máxima"
$ws.Cells.Item(68, 2).Value = 3
$ws.Cells.Item(69, 1).Value = "Demanda total"
$ws.Cells.Item(69, 2).Value = 6
$ws.Cells.Item(70, 1).Value = "Duración Media (Dm)"
$ws.Cells.Item(70, 2).Value = "0:00"
